$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings need to be forced to
# Text format first, otherwise Excel auto-converts them to numbers and
# loses the exact original text representation (e.g. trailing zeros).
$ws.Range("D2").Value = '43.139.17'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '2.348.88'
$ws.Range("E3").Value = '  +2.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.65'
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.38'
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.13'
$ws.Range("E10").Value = '  -1.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0787'
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.57'
$ws.Range("E12").Value = '  -3.07%  '
$ws.Range("E13").Value = '  +2.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.75'
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").Value = '2.715.13'
$ws.Range("E15").Value = '  +2.41%  '
$ws.Range("D16").Value = '2.327.37'
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.797'
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("D18").Value = '43.102.94'
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.18'
$ws.Range("E19").Value = '  -0.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.22'
$ws.Range("E20").Value = '  +3.30%  '
$ws.Range("D21").Value = '0.0₃0889'
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.09'
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.43'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  -1.73%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.59'
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  +14.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.16'
$ws.Range("E29").Value = '  +1.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.25'
$ws.Range("E30").Value = '  -3.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0722'
$ws.Range("E33").Value = '  +3.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.19'
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("E35").Value = '  +5.04%  '
$ws.Range("E36").Value = '  -1.80%  '
$ws.Range("E38").Value = '  +0.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.66'
$ws.Range("E39").Value = '  +14.52%  '
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '111.21'
$ws.Range("E42").Value = '  -32.44%  '
$ws.Range("D43").Value = '1.936.21'
$ws.Range("E43").Value = '  -1.10%  '
$ws.Range("E44").Value = '  +0.60%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.11'
$ws.Range("E45").Value = '  +3.25%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.92'
$ws.Range("E46").Value = '  -5.58%  '
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("D48").Value = '2.578.41'
$ws.Range("E48").Value = '  +2.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '52.88'
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.78'
$ws.Range("E50").Value = '  -3.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.21'
$ws.Range("E51").Value = '  +0.71%  '
